$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''43.011.59'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.94%  '

$ws.Range("D3").Value = '''2.339.40'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.15%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").Value = '''306.77'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.51%  '

$ws.Range("D6").Value = '''101.00'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.97%  '

$ws.Range("E7").Value = '  -4.10%  '

$ws.Range("E8").Value = '  -0.01%  '

$ws.Range("E9").Value = '  -3.62%  '

$ws.Range("D10").Value = '''34.92'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.31%  '

$ws.Range("D11").Value = '''52.49'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.83%  '

$ws.Range("E12").Value = '  -2.09%  '

$ws.Range("E13").Value = '  +0.73%  '

$ws.Range("D14").Value = '''6.87'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.42%  '

$ws.Range("D15").Value = '''15.82'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +5.20%  '

$ws.Range("D16").Value = '''2.333.95'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.46%  '

$ws.Range("D17").Value = '''0.829'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.28%  '

$ws.Range("D18").Value = '''42.937.48'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.90%  '

$ws.Range("E19").Value = '  +0.79%  '

$ws.Range("D20").Value = '''11.75'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.30%  '

$ws.Range("D21").Value = '''0.0₃0911'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.48%  '

$ws.Range("D22").Value = '''68.02'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.20%  '

$ws.Range("D23").Value = '''236.70'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.08%  '

$ws.Range("D24").Value = '''2.03'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.81%  '

$ws.Range("D25").Value = '''2.56'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.92%  '

$ws.Range("E26").Value = '  -0.14%  '

$ws.Range("D27").Value = '''25.37'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.62%  '

$ws.Range("D28").Value = '''3.96'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.44%  '

$ws.Range("E29").Value = '  +1.02%  '

$ws.Range("D30").Value = '''35.46'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.58%  '

$ws.Range("D31").Value = '''9.32'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.52%  '

$ws.Range("D32").Value = '''164.01'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.96%  '

$ws.Range("E33").Value = '  -0.05%  '

$ws.Range("B35").Value = 'RenderToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D35").Value = '''4.66'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +7.79%  '

$ws.Range("B36").Value = 'Celestia'
$ws.Range("C36").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D36").Value = '''17.60'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.05%  '

$ws.Range("D37").Value = '''0.0728'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.80%  '

$ws.Range("E38").Value = '  -4.37%  '

$ws.Range("E39").Value = '  -1.54%  '

$ws.Range("E40").Value = '  -5.10%  '

$ws.Range("E41").Value = '  -3.08%  '

$ws.Range("E42").Value = '  -2.55%  '

$ws.Range("D43").Value = '''2.51'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +9.11%  '

$ws.Range("D44").Value = '''2.029.51'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.67%  '

$ws.Range("E45").Value = '  -2.31%  '

$ws.Range("D46").Value = '''18.92'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.06%  '

$ws.Range("D47").Value = '''10.19'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.68%  '

$ws.Range("E48").Value = '  -2.10%  '

$ws.Range("D49").Value = '''56.60'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.94%  '

$ws.Range("E50").Value = '  -1.39%  '

$ws.Range("D51").Value = '''2.564.03'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.06%  '

